# Update "想去人数" (want-to-go count) figures in column F across the
# three sheets that carry them (展览, 演出, 全部类型) to match the
# refreshed stats baked into this gh-pages build.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (index 1) ---
$ws1 = $wb.Worksheets.Item(1)
$sheet1Updates = @{
    4  = 444
    5  = 1304
    7  = 7626
    8  = 93
    10 = 2087
    11 = 8255
    12 = 8
    15 = 5640
    17 = 2598
    21 = 400
    24 = 521
    25 = 3450
    26 = 38
    28 = 18
    29 = 2935
    30 = 10
    31 = 55
    33 = 337
    34 = 127
    35 = 298
    36 = 170
    37 = 651
    40 = 1663
    44 = 2704
    45 = 5
}
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Range("F$row").Value = $sheet1Updates[$row]
}

# --- Sheet "演出" (index 2) ---
$ws2 = $wb.Worksheets.Item(2)
$sheet2Updates = @{
    3 = 117
    4 = 5
    5 = 46
}
foreach ($row in $sheet2Updates.Keys) {
    $ws2.Range("F$row").Value = $sheet2Updates[$row]
}

# --- Sheet "全部类型" (index 4) ---
$ws4 = $wb.Worksheets.Item(4)
$sheet4Updates = @{
    6  = 1304
    7  = 7626
    8  = 93
    10 = 2087
    11 = 8256
    12 = 8
    15 = 5640
    17 = 2598
    20 = 400
    24 = 117
    25 = 521
    26 = 3451
    28 = 18
    29 = 2935
    30 = 337
    31 = 127
    32 = 298
    33 = 46
    34 = 651
    39 = 1663
    43 = 2704
    44 = 5
}
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Range("F$row").Value = $sheet4Updates[$row]
}
